$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.031528711567798
$ws.Range("D2").Value = 1.040777922395563
$ws.Range("E2").Value = 1.031004210401216
$ws.Range("F2").Value = 1.047580054212719
$ws.Range("J2").Value = 1.036663938705475
$ws.Range("K2").Value = 1.043559307405299
$ws.Range("L2").Value = 1.033813583220028
$ws.Range("M2").Value = 1.050342295136893
$ws.Range("N2").Value = 1.038136120405452
# Row 3
$ws.Range("C3").Value = 1.033149473229332
$ws.Range("D3").Value = 1.042373664865803
$ws.Range("E3").Value = 1.032405029692096
$ws.Range("F3").Value = 1.049331019988978
$ws.Range("J3").Value = 1.037923198507259
$ws.Range("K3").Value = 1.04496347961889
$ws.Range("L3").Value = 1.03502128702245
$ws.Range("M3").Value = 1.051902701189665
$ws.Range("N3").Value = 1.039397168500594
# Row 4
$ws.Range("C4").Value = 1.034196746797846
$ws.Range("D4").Value = 1.04340511938952
$ws.Range("E4").Value = 1.033310310487998
$ws.Range("F4").Value = 1.05046310950982
$ws.Range("J4").Value = 1.038736259236987
$ws.Range("K4").Value = 1.045870500420915
$ws.Range("L4").Value = 1.035801096219005
$ws.Range("M4").Value = 1.052911040437299
$ws.Range("N4").Value = 1.040211383869818
# Row 5
$ws.Range("C5").Value = 1.034636679974208
$ws.Range("D5").Value = 1.043838489932027
$ws.Range("E5").Value = 1.033690625064756
$ws.Range("F5").Value = 1.050938834856217
$ws.Range("J5").Value = 1.039077655643758
$ws.Range("K5").Value = 1.046251444461984
$ws.Range("L5").Value = 1.036128539092528
$ws.Range("M5").Value = 1.053334633519881
$ws.Range("N5").Value = 1.040553265098645
# Row 6
$ws.Range("C6").Value = 1.034710527025687
$ws.Range("D6").Value = 1.04391124020476
$ws.Range("E6").Value = 1.033754466224417
$ws.Range("F6").Value = 1.051018699525774
$ws.Range("J6").Value = 1.039134953541414
$ws.Range("K6").Value = 1.0463153853514
$ws.Range("L6").Value = 1.03618349560334
$ws.Range("M6").Value = 1.053405738610759
$ws.Range("N6").Value = 1.040610644365888
# Row 7
$ws.Range("C7").Value = 1.034202626530107
$ws.Range("D7").Value = 1.04341091108818
$ws.Range("E7").Value = 1.033315393301748
$ws.Range("F7").Value = 1.050469466965876
$ws.Range("J7").Value = 1.038740822609221
$ws.Range("K7").Value = 1.045875592045446
$ws.Range("L7").Value = 1.03580547304441
$ws.Range("M7").Value = 1.052916701723882
$ws.Range("N7").Value = 1.040215953722564
# Row 8
$ws.Range("C8").Value = 1.032076765904279
$ws.Range("D8").Value = 1.041317444217921
$ws.Range("E8").Value = 1.03147786541934
$ws.Range("F8").Value = 1.048171994757701
$ws.Range("J8").Value = 1.037089881716239
$ws.Range("K8").Value = 1.044034184986658
$ws.Range("L8").Value = 1.03422207978242
$ws.Range("M8").Value = 1.05086992673562
$ws.Range("N8").Value = 1.038562668304149
# Row 9
$ws.Range("C9").Value = 1.028319006331302
$ws.Range("D9").Value = 1.037619640883565
$ws.Range("E9").Value = 1.028230786132571
$ws.Range("F9").Value = 1.044116151614282
$ws.Range("J9").Value = 1.034166833571059
$ws.Range("K9").Value = 1.04077694422093
$ws.Range("L9").Value = 1.031418918133993
$ws.Range("M9").Value = 1.0472524788012
$ws.Range("N9").Value = 1.035635469095323
# Row 10
$ws.Range("C10").Value = 1.025805277719833
$ws.Range("D10").Value = 1.035147871614436
$ws.Range("E10").Value = 1.026059419911653
$ws.Range("F10").Value = 1.041406588782684
$ws.Range("J10").Value = 1.032208291231636
$ws.Range("K10").Value = 1.038596517508378
$ws.Range("L10").Value = 1.029540915074887
$ws.Range("M10").Value = 1.044832980821994
$ws.Range("N10").Value = 1.03367414540113
# Row 11
$ws.Range("C11").Value = 1.024714625673881
$ws.Range("D11").Value = 1.034075872470286
$ws.Range("E11").Value = 1.025117502779489
$ws.Range("F11").Value = 1.04023182003385
$ws.Range("J11").Value = 1.031357771619474
$ws.Range("K11").Value = 1.037650123774023
$ws.Range("L11").Value = 1.028725424021697
$ws.Range("M11").Value = 1.043783306098854
$ws.Range("N11").Value = 1.032822417953563
# Row 12
$ws.Range("C12").Value = 1.024309166753347
$ws.Range("D12").Value = 1.033677415882243
$ws.Range("E12").Value = 1.02476736724914
$ws.Range("F12").Value = 1.039795218422662
$ws.Range("J12").Value = 1.031041471582675
$ws.Range("K12").Value = 1.037298241350387
$ws.Range("L12").Value = 1.028422159079493
$ws.Range("M12").Value = 1.043393095180532
$ws.Range("N12").Value = 1.032505668734432
# Row 13
$ws.Range("C13").Value = 1.024396154763193
$ws.Range("D13").Value = 1.033762898510581
$ws.Range("E13").Value = 1.02484248468764
$ws.Range("F13").Value = 1.039888882034356
$ws.Range("J13").Value = 1.031109336333286
$ws.Range("K13").Value = 1.037373737323711
$ws.Range("L13").Value = 1.028487226663756
$ws.Range("M13").Value = 1.043476811196233
$ws.Range("N13").Value = 1.032573629860773
# Row 14
$ws.Range("C14").Value = 1.024681117369311
$ws.Range("D14").Value = 1.034042941468731
$ws.Range("E14").Value = 1.025088565942883
$ws.Range("F14").Value = 1.040195735385686
$ws.Range("J14").Value = 1.031331633952679
$ws.Range("K14").Value = 1.037621044259477
$ws.Range("L14").Value = 1.028700363347615
$ws.Range("M14").Value = 1.043751057638611
$ws.Range("N14").Value = 1.032796243168284
# Row 15
$ws.Range("C15").Value = 1.024856646501032
$ws.Range("D15").Value = 1.034215449256044
$ws.Range("E15").Value = 1.025240149235509
$ws.Range("F15").Value = 1.040384765635167
$ws.Range("J15").Value = 1.031468548309035
$ws.Range("K15").Value = 1.037773371563875
$ws.Range("L15").Value = 1.028831636541384
$ws.Range("M15").Value = 1.043919987820463
$ws.Range("N15").Value = 1.03293335195873
# Row 16
$ws.Range("C16").Value = 1.025877613434062
$ws.Range("D16").Value = 1.035218979689826
$ws.Range("E16").Value = 1.02612189515696
$ws.Range("F16").Value = 1.041484521364972
$ws.Range("J16").Value = 1.03226468478837
$ws.Range("K16").Value = 1.038659278111699
$ws.Range("L16").Value = 1.029594987217837
$ws.Range("M16").Value = 1.044902600769329
$ws.Range("N16").Value = 1.033730619043184
# Row 17
$ws.Range("C17").Value = 1.02651744327494
$ws.Range("D17").Value = 1.035848002466259
$ws.Range("E17").Value = 1.026674528593647
$ws.Range("F17").Value = 1.04217395540717
$ws.Range("J17").Value = 1.03276341515529
$ws.Range("K17").Value = 1.039214372583785
$ws.Range("L17").Value = 1.030073193787925
$ws.Range("M17").Value = 1.045518419634303
$ws.Range("N17").Value = 1.034230057664429
# Row 18
$ws.Range("C18").Value = 1.026890434684554
$ws.Range("D18").Value = 1.036214736947007
$ws.Range("E18").Value = 1.026996706862235
$ws.Range("F18").Value = 1.042575945815628
$ws.Range("J18").Value = 1.033054079667692
$ws.Range("K18").Value = 1.039537933006811
$ws.Range("L18").Value = 1.030351901977319
$ws.Range("M18").Value = 1.045877422653877
$ws.Range("N18").Value = 1.034521134953779
# Row 19
$ws.Range("C19").Value = 1.02701757980203
$ws.Range("D19").Value = 1.036339756493795
$ws.Range("E19").Value = 1.027106533899832
$ws.Range("F19").Value = 1.042712989996006
$ws.Range("J19").Value = 1.033153148943727
$ws.Range("K19").Value = 1.039648222307022
$ws.Range("L19").Value = 1.030446896955261
$ws.Range("M19").Value = 1.045999800914292
$ws.Range("N19").Value = 1.03462034491955
# Row 20
$ws.Range("C20").Value = 1.026448817456672
$ws.Range("D20").Value = 1.035780531310927
$ws.Range("E20").Value = 1.026615253236135
$ws.Range("F20").Value = 1.042100000733194
$ws.Range("J20").Value = 1.032709930638099
$ws.Range("K20").Value = 1.039154838703062
$ws.Range("L20").Value = 1.03002190973766
$ws.Range("M20").Value = 1.045452368221459
$ws.Range("N20").Value = 1.03417649719309
# Row 21
$ws.Range("C21").Value = 1.024597212535539
$ws.Range("D21").Value = 1.033960483302447
$ws.Range("E21").Value = 1.025016108558191
$ws.Range("F21").Value = 1.040105381459722
$ws.Range("J21").Value = 1.031266183391792
$ws.Range("K21").Value = 1.037548228271203
$ws.Range("L21").Value = 1.028637609787623
$ws.Range("M21").Value = 1.043670307689586
$ws.Range("N21").Value = 1.032730699660094
# Row 22
$ws.Range("C22").Value = 1.023431048698817
$ws.Range("D22").Value = 1.032814588148412
$ws.Range("E22").Value = 1.024009121324592
$ws.Range("F22").Value = 1.03884988880944
$ws.Range("J22").Value = 1.030356243381403
$ws.Range("K22").Value = 1.036536060063308
$ws.Range("L22").Value = 1.027765185449951
$ws.Range("M22").Value = 1.042548027008647
$ws.Range("N22").Value = 1.031819467430518
# Row 23
$ws.Range("C23").Value = 1.024049446870192
$ws.Range("D23").Value = 1.033422200618802
$ws.Range("E23").Value = 1.024543093780496
$ws.Range("F23").Value = 1.039515586047211
$ws.Range("J23").Value = 1.030838831553191
$ws.Range("K23").Value = 1.037072825528397
$ws.Range("L23").Value = 1.028227872457807
$ws.Range("M23").Value = 1.043143146614359
$ws.Range("N23").Value = 1.032302740932863
# Row 24
$ws.Range("C24").Value = 1.026479827160455
$ws.Range("D24").Value = 1.035811019128139
$ws.Range("E24").Value = 1.026642037723909
$ws.Range("F24").Value = 1.04213341811449
$ws.Range("J24").Value = 1.032734098719935
$ws.Range("K24").Value = 1.039181740171069
$ws.Range("L24").Value = 1.030045083477312
$ws.Range("M24").Value = 1.045482214610346
$ws.Range("N24").Value = 1.034200699596374
# Row 25
$ws.Range("C25").Value = 1.029291935457606
$ws.Range("D25").Value = 1.038576723521798
$ws.Range("E25").Value = 1.029071367504201
$ws.Range("F25").Value = 1.045165632324112
$ws.Range("J25").Value = 1.034924208806535
$ws.Range("K25").Value = 1.041620551846369
$ws.Range("L25").Value = 1.032145194206915
$ws.Range("M25").Value = 1.048189016084253
$ws.Range("N25").Value = 1.036393919890506

Write-Host "Applied 380 kV case values"